# Apply leaderboard/building-pool board updates (#21, #18)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: HSE -> MON block, FAC -> MON blocks
$ws.Range("D3").Value = "n"
$ws.Range("E3").Value = "MON"
$ws.Range("F3").Value = "n"
$ws.Range("K3").Value = "MON"
$ws.Range("Q3").Value = "MON"
$ws.Range("W3").Value = "MON"

# Row 5: HWY block (was SHP / HSE)
$ws.Range("D5").Value = "n"
$ws.Range("E5").Value = "HWY"
$ws.Range("F5").Value = "n"
$ws.Range("Q5").Value = "HWY"
$ws.Range("W5").Value = "HWY"

# Row 7: PRK block (was BCH / HSE / FAC)
$ws.Range("D7").Value = "n"
$ws.Range("E7").Value = "PRK"
$ws.Range("F7").Value = "n"
$ws.Range("K7").Value = "PRK"
$ws.Range("Q7").Value = "PRK"
$ws.Range("W7").Value = "PRK"

# Row 9: MON block newly filled in (was empty "e"), PRK block newly filled in, BCH -> MON
$ws.Range("D9").Value = "n"
$ws.Range("E9").Value = "MON"
$ws.Range("F9").Value = "n"
$ws.Range("J9").Value = "n"
$ws.Range("K9").Value = "PRK"
$ws.Range("L9").Value = "n"
$ws.Range("W9").Value = "MON"
